# Fixed bug for RestDay penalization
# Update TabuSearch_Stats worksheet with corrected optimization statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabu Stats")

# ---- Table 1: per-division stats (rows 2-12) ----
# Row 2: Division U7 (tier: 1)
$ws.Range("B2").Value = 205.86
$ws.Range("C2").Value = 10517.0
$ws.Range("D2").Value = 48.0
$ws.Range("E2").Value = 10469.0
$ws.Range("F2").Value = "15 sec"

# Row 3: Division U7 (tier: 2)
$ws.Range("B3").Value = 221.82
$ws.Range("C3").Value = 792.0
$ws.Range("D3").Value = 30.0
$ws.Range("E3").Value = 762.0
$ws.Range("F3").Value = "4 sec"

# Row 4: Division U7 (tier: 3)
$ws.Range("B4").Value = 212.33999999999997
$ws.Range("C4").Value = 1002.0
$ws.Range("D4").Value = 39.0
$ws.Range("E4").Value = 963.0
$ws.Range("F4").Value = "4 sec"

# Row 5: Division U8 (tier: 0)
$ws.Range("B5").Value = 334.45
$ws.Range("C5").Value = 5735.0
$ws.Range("E5").Value = 5731.0
$ws.Range("F5").Value = "12 sec"

# Row 6: Division U8 (tier: 1)
$ws.Range("B6").Value = 386.99
$ws.Range("C6").Value = 11965.0
$ws.Range("E6").Value = 11965.0
$ws.Range("F6").Value = "12 sec"

# Row 7: Division U8 (tier: 2)
$ws.Range("B7").Value = 291.37
$ws.Range("C7").Value = 14962.0
$ws.Range("D7").Value = 30.0
$ws.Range("E7").Value = 14932.0
$ws.Range("F7").Value = "14 sec"

# Row 8: Division U8 (tier: 3)
$ws.Range("B8").Value = 175.26999999999998
$ws.Range("C8").Value = 10993.0
$ws.Range("D8").Value = 33.0
$ws.Range("E8").Value = 10960.0
$ws.Range("F8").Value = "15 sec"

# Row 9: Division U9 (tier: 0)
$ws.Range("B9").Value = 181.43
$ws.Range("C9").Value = 966.0
$ws.Range("E9").Value = 966.0
$ws.Range("F9").Value = "5 sec"

# Row 10: Division U9 (tier: 1)
$ws.Range("B10").Value = 203.73000000000002
$ws.Range("C10").Value = 4599.0
$ws.Range("D10").Value = 44.0
$ws.Range("E10").Value = 4555.0
$ws.Range("F10").Value = "12 sec"

# Row 11: Division U9 (tier: 2)
$ws.Range("B11").Value = 230.19
$ws.Range("C11").Value = 1332.0
$ws.Range("D11").Value = 32.0
$ws.Range("E11").Value = 1300.0
$ws.Range("F11").Value = "6 sec"

# Row 12: Division U9 (tier: 3)
$ws.Range("B12").Value = 435.58
$ws.Range("C12").Value = 12388.0
$ws.Range("D12").Value = 201.0
$ws.Range("E12").Value = 12187.0
$ws.Range("F12").Value = "14 sec"

# ---- Table 2: per-division stats, second schedule (rows 16-26) ----
# Row 16: Division U7 (tier: 1)
$ws.Range("B16").Value = 134.66
$ws.Range("C16").Value = 1554.0
$ws.Range("D16").Value = 68.0
$ws.Range("E16").Value = 1486.0

# Row 17: Division U7 (tier: 2)
$ws.Range("B17").Value = 152.03
$ws.Range("C17").Value = 1079.0
$ws.Range("D17").Value = 86.0
$ws.Range("E17").Value = 993.0

# Row 18: Division U7 (tier: 3)
$ws.Range("B18").Value = 176.05
$ws.Range("C18").Value = 3145.0
$ws.Range("D18").Value = 141.0
$ws.Range("E18").Value = 3004.0

# Row 19: Division U8 (tier: 0)
$ws.Range("B19").Value = 177.37
$ws.Range("C19").Value = 957.0
$ws.Range("D19").Value = 40.0
$ws.Range("E19").Value = 917.0

# Row 20: Division U8 (tier: 1)
$ws.Range("B20").Value = 225.81
$ws.Range("C20").Value = 537.0
$ws.Range("D20").Value = 46.0
$ws.Range("E20").Value = 491.0

# Row 21: Division U8 (tier: 2)
$ws.Range("B21").Value = 167.01999999999998
$ws.Range("C21").Value = 1184.0
$ws.Range("D21").Value = 82.0
$ws.Range("E21").Value = 1102.0

# Row 22: Division U8 (tier: 3)
$ws.Range("B22").Value = 149.96
$ws.Range("C22").Value = 1342.0
$ws.Range("D22").Value = 77.0
$ws.Range("E22").Value = 1265.0

# Row 23: Division U9 (tier: 0)
$ws.Range("B23").Value = 131.89999999999998
$ws.Range("C23").Value = 885.0
$ws.Range("D23").Value = 44.0
$ws.Range("E23").Value = 841.0

# Row 24: Division U9 (tier: 1)
$ws.Range("B24").Value = 155.0
$ws.Range("C24").Value = 1202.0
$ws.Range("D24").Value = 102.0
$ws.Range("E24").Value = 1100.0

# Row 25: Division U9 (tier: 2)
$ws.Range("B25").Value = 172.26999999999998
$ws.Range("C25").Value = 1298.0
$ws.Range("D25").Value = 75.0
$ws.Range("E25").Value = 1223.0

# Row 26: Division U9 (tier: 3)
$ws.Range("B26").Value = 230.3
$ws.Range("C26").Value = 4046.0
$ws.Range("D26").Value = 213.0
$ws.Range("E26").Value = 3833.0

# ---- Row 30: Entire League summary ----
$ws.Range("B30").Value = 423.19
$ws.Range("C30").Value = 17229.0
$ws.Range("D30").Value = 974.0
$ws.Range("E30").Value = 16255.0
$ws.Range("F30").Value = "2 min, 13 sec"
